$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '56.813.84'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.05%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.967.56'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -1.19%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '497.49'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -3.06%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '136.78'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -1.46%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.20%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -2.35%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '7.30'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -3.16%  '
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -1.85%  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -0.40%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '3.475.94'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -1.17%  '
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -1.61%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '25.73'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.29%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000157'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -0.25%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '56.920.01'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.10%  '
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +1.80%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.956.81'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -1.60%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.57'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.12%  '
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -0.90%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '318.04'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -2.76%  '
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.24%  '
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.82%  '
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.26%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '62.92'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.70%  '
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.21%  '
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -5.32%  '
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -3.66%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.52'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -1.79%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.06'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -0.37%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.76'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -2.73%  '
$ws.Range("B32").NumberFormat = "@"
$ws.Range("B32").Value = 'Fetch.AI'
$ws.Range("C32").NumberFormat = "@"
$ws.Range("C32").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.15'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -6.66%  '
$ws.Range("B33").NumberFormat = "@"
$ws.Range("B33").Value = 'EthereumClassic'
$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '20.11'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -2.37%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '154.45'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -1.96%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.60'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.44%  '
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +0.01%  '
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -2.30%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '23.90'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -0.82%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0663'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -2.32%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.000.17'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -1.17%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '37.41'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +0.69%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.00'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +0.07%  '
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +0.66%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.636'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -2.02%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.192.48'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -4.14%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.38'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -3.37%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.938'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -6.27%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '5.92'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +0.73%  '
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -2.58%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '19.09'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -0.53%  '
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -10.06%  '
